$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "27.029.40"
$ws.Range("E2").Value = "  -2.08%  "
Set-TextValue $ws.Range("D3") "1.824.63"
$ws.Range("E3").Value = "  -0.97%  "
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  -0.74%  "
Set-TextValue $ws.Range("D5") "311.47"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("E6").Value = "  -0.65%  "
Set-TextValue $ws.Range("D7") "0.4245"
$ws.Range("E7").Value = "  -1.44%  "
Set-TextValue $ws.Range("D8") "0.3667"
$ws.Range("E8").Value = "  -1.63%  "
Set-TextValue $ws.Range("D9") "0.07227"
$ws.Range("E9").Value = "  -0.87%  "
Set-TextValue $ws.Range("D10") "0.8422"
$ws.Range("E10").Value = "  -3.29%  "
Set-TextValue $ws.Range("D11") "20.71"
$ws.Range("E11").Value = "  -2.82%  "
Set-TextValue $ws.Range("D12") "1.834.96"
$ws.Range("E12").Value = "  -0.54%  "
Set-TextValue $ws.Range("D13") "6.665"
$ws.Range("E13").Value = "  -0.79%  "
Set-TextValue $ws.Range("D14") "5.284"
$ws.Range("E14").Value = "  -1.88%  "
Set-TextValue $ws.Range("D15") "0.07036"
$ws.Range("E15").Value = "  -1.30%  "
Set-TextValue $ws.Range("D16") "89.63"
$ws.Range("E16").Value = "  +1.19%  "
Set-TextValue $ws.Range("D17") "1.001"
$ws.Range("E17").Value = "  -0.90%  "
Set-TextValue $ws.Range("D18") "0.000008741"
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("E19").Value = "  -0.66%  "
Set-TextValue $ws.Range("D20") "14.88"
$ws.Range("E20").Value = "  -2.89%  "
Set-TextValue $ws.Range("D21") "27.131.54"
$ws.Range("E21").Value = "  -1.71%  "
Set-TextValue $ws.Range("D22") "5.127"
$ws.Range("E22").Value = "  -1.03%  "
Set-TextValue $ws.Range("D23") "10.80"
$ws.Range("E23").Value = "  -1.84%  "
Set-TextValue $ws.Range("D24") "2.054.29"
$ws.Range("E24").Value = "  -0.73%  "
Set-TextValue $ws.Range("D25") "1.979"
$ws.Range("E25").Value = "  +0.70%  "
Set-TextValue $ws.Range("D26") "151.37"
$ws.Range("E26").Value = "  -1.92%  "
Set-TextValue $ws.Range("D27") "2.247"
$ws.Range("E27").Value = "  +4.18%  "
Set-TextValue $ws.Range("D28") "18.17"
$ws.Range("E28").Value = "  -1.92%  "
Set-TextValue $ws.Range("D29") "5.235"
$ws.Range("E29").Value = "  -1.38%  "
Set-TextValue $ws.Range("D30") "116.76"
$ws.Range("E30").Value = "  -0.65%  "
Set-TextValue $ws.Range("D31") "0.08716"
$ws.Range("E31").Value = "  -2.06%  "
Set-TextValue $ws.Range("D32") "1.179"
$ws.Range("E32").Value = "  -2.75%  "
Set-TextValue $ws.Range("D33") "0.7353"
$ws.Range("E33").Value = "  -4.73%  "
$ws.Range("E34").Value = "  -0.18%  "
Set-TextValue $ws.Range("D35") "4.419"
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("E36").Value = "  -0.82%  "
Set-TextValue $ws.Range("D37") "1.093"
$ws.Range("E37").Value = "  -3.00%  "
$ws.Range("E38").Value = "  -1.17%  "
Set-TextValue $ws.Range("D39") "0.05227"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("E40").Value = "  +2.16%  "
Set-TextValue $ws.Range("D41") "2.870"
$ws.Range("E41").Value = "  -0.55%  "
Set-TextValue $ws.Range("D42") "0.5131"
$ws.Range("E42").Value = "  +0.47%  "
Set-TextValue $ws.Range("D43") "0.1686"
$ws.Range("E43").Value = "  +0.08%  "
Set-TextValue $ws.Range("D44") "8.545"
$ws.Range("E44").Value = "  -2.29%  "
Set-TextValue $ws.Range("D45") "10.54"
$ws.Range("E45").Value = "  -0.75%  "
Set-TextValue $ws.Range("D46") "1.954"
$ws.Range("E46").Value = "  +6.33%  "
Set-TextValue $ws.Range("D47") "0.4732"
$ws.Range("E47").Value = "  -0.10%  "
Set-TextValue $ws.Range("D48") "105.59"
$ws.Range("E48").Value = "  -1.20%  "
Set-TextValue $ws.Range("D49") "1.000"
$ws.Range("E49").Value = "  -0.80%  "
Set-TextValue $ws.Range("D50") "0.06325"
$ws.Range("E50").Value = "  -1.83%  "
Set-TextValue $ws.Range("D51") "1.650"
$ws.Range("E51").Value = "  -1.66%  "
